$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") rows 2-38, replacing old "Strike#"-derived
# values with the regenerated K values.
$gValues = @{
    2  = 2
    3  = 1
    4  = 2
    5  = 8
    6  = 2
    7  = 5
    8  = 7
    9  = 6
    10 = 3
    11 = 2
    12 = 2
    13 = 5
    14 = 6
    15 = 6
    16 = 2
    17 = 9
    18 = 3
    19 = 3
    20 = 5
    21 = 4
    22 = 2
    23 = 6
    24 = 3
    25 = 1
    26 = 3
    27 = 6
    28 = 2
    29 = 6
    30 = 5
    31 = 3
    32 = 3
    33 = 3
    34 = 2
    35 = 5
    36 = 3
    37 = 3
    38 = 2
}

foreach ($row in $gValues.Keys) {
    $ws.Range("G$row").Value = $gValues[$row]
}
